$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.344.58'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9978'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.91'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6283'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9992'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07445'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2904'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.42'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07730'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.985'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6787'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.03'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.186'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.378.54'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.74'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.32'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.501'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9991'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.82'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.495'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1368'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06442'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +14.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.416'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.480'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.088'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.837'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.140'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6942'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.578'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.258.42'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.835'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01837'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.777'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9316'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9986'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.009.79'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.82'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.47%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.723'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.062'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1156'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.991'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3935'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.59%  '
